$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17-18: WrappedEther and ShibaInu swap positions with updated data
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.11%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.668.53"
$ws.Range("E18").Value = "  +2.01%  "

# Remaining price / volume updates
$ws.Range("D2").Value = "60.855.33"
$ws.Range("E2").Value = "  +4.02%  "
$ws.Range("D3").Value = "2.645.70"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.55%  "
$ws.Range("D9").Value = "2.671.01"
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.86"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("E11").Value = "  +5.57%  "
$ws.Range("E12").Value = "  +6.92%  "
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").Value = "3.117.21"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").Value = "60.494.21"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.15%  "
$ws.Range("E19").Value = "  +3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.06"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.49"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.81%  "
$ws.Range("E22").Value = "  +2.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.74"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.442"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.93%  "
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.990"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  +5.07%  "
$ws.Range("D30").Value = "0.0₃0815"
$ws.Range("E30").Value = "  +11.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  +4.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.26"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.54"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("E36").Value = "  +5.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.908"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.909"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +11.87%  "
$ws.Range("E39").Value = "  +5.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.60"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("E41").Value = "  +7.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "303.92"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.79%  "
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0985"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.608"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("E47").Value = "  +4.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.82"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +14.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  +5.40%  "
